$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B20").Value = 44835
$ws.Range("C20").Value = 0.90625
$ws.Range("D20").Value = 0.90625
$ws.Range("D20").Font.Bold = $true
$ws.Range("D20").Font.Italic = $true
$ws.Range("D20").NumberFormat = "h:mm AM/PM"
$ws.Range("E20").Formula = "=D20-C20"
